# Latest data import scripts + sheets.
# Rebuild the Room Legend table (columns A:B) to match the refreshed room
# data export, re-style the two-section legend tail (rows 50-56) with a
# slightly larger font, fix up the hidden AutoFilter defined name, and
# reset stray custom row heights left over from the old layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Full target data for A1:B56 (room code -> room description).
# ---------------------------------------------------------------------
$rows = @(
    @(1, "Y-CR1", "SAGA COLLEGE - CLASSROOM 1"),
    @(2, "Y-CR2", "SAGA COLLEGE - CLASSROOM 2"),
    @(3, "Y-CR3", "SAGA COLLEGE - CLASSROOM 3"),
    @(4, "Y-CR4", "SAGA COLLEGE - CLASSROOM 4"),
    @(5, "Y-CR5", "SAGA COLLEGE - CLASSROOM 5"),
    @(6, "Y-CR6", "SAGA COLLEGE - CLASSROOM 6"),
    @(7, "Y-CR7", "SAGA COLLEGE - CLASSROOM 7"),
    @(8, "Y-CR8", "SAGA COLLEGE - CLASSROOM 8"),
    @(9, "Y-CR9", "ELM COLLEGE - CLASSROOM 9"),
    @(10, "Y-KChanrai", "ELM COLLEGE - KEWALRAM CHANRAI ROOM"),
    @(11, "Y-CR11", "ELM COLLEGE - CLASSROOM 11"),
    @(12, "Y-CR12", "ELM COLLEGE - CLASSROOM 12"),
    @(13, "Y-CR13", "ELM COLLEGE - CLASSROOM 13"),
    @(14, "Y-CR14", "ELM COLLEGE - CLASSROOM 14"),
    @(15, "Y-CR15", "LIBRARY - CLASSROOM 15"),
    @(16, "Y-CR16", "LIBRARY - CLASSROOM 16"),
    @(17, "Y-CR17", "SCIENCE CENTRE - CLASSROOM 17"),
    @(18, "Y-CR18", "CENDANA COLLEGE - CLASSROOM 18"),
    @(19, "Y-CR19", "CENDANA COLLEGE - CLASSROOM 19"),
    @(20, "Y-CR20", "CENDANA COLLEGE - CLASSROOM 20"),
    @(21, "Y-CR21", "CENDANA COLLEGE - CLASSROOM 21"),
    @(22, "Y-CR22", "CENDANA COLLEGE - CLASSROOM 22"),
    @(23, "Y-CR23", "CENDANA COLLEGE - CLASSROOM 23"),
    @(24, "Y-GLRm1", "LIBRARY - GLOBAL LEARNING ROOM 1"),
    @(25, "Y-GLRm2", "LIBRARY - GLOBAL LEARNING ROOM 2"),
    @(26, "Y-PgRm1", "LIBRARY - PROGRAMME ROOM 1"),
    @(27, "Y-PgRm2", "LIBRARY - PROGRAMME ROOM 2"),
    @(28, "Y-CompLab", "LIBRARY - COMPUTER LAB"),
    @(29, "Y-LT1", "SAGA COLLEGE - LECTURE THEATRE 1"),
    @(30, "Y-TCTLT", "ELM COLLEGE - TAN CHIN TUAN LECTURE THEATRE"),
    @(31, "Y-BioLab", "SCIENCE CENTRE - LAB 1 (Physics Lab)"),
    @(32, "Y-AChemLab", "SCIENCE CENTRE - LAB 2 (Life Sciences Lab)"),
    @(33, "Y-PhysLab", "SCIENCE CENTRE - LAB 3 (Y.E.S. Lab)"),
    @(34, "Y-OChemLab", "SCIENCE CENTRE - LAB 4 (Chemistry Lab)"),
    @(35, "Y-ArtsStud", "CENDANA COLLEGE - STUDIO 2 (GEORGETTE CHEN ARTS STUDIO)"),
    @(36, "Y-Studio3", "CENDANA COLLEGE - STUDIO 3"),
    @(37, "Y-DanceStu", "CENDANA COLLEGE - STUDIO 1 (DANCE STUDIO)"),
    @(38, "Y-Studio4", "ARTS CENTRE - FABRICATION STUDIO 4"),
    @(39, "Y-Studio5", "ARTS CENTRE - FABRICATION STUDIO 5"),
    @(40, "Y-PracRm6", "ARTS CENTRE - PRACTICE ROOM 6"),
    @(41, "Y-PerfHall", "ARTS CENTRE - PERFORMANCE HALL"),
    @(42, "Y-RC2SC", "RC2 Student Common"),
    @(43, "UT-AUD03", "UT- Auditorium 3"),
    @(44, "UTSRC-AUD2", "SRC - Auditorium 2"),
    @(45, "UTSRC-GLR", "SRC - Global Learning Room"),
    @(46, "UTSRC-LT51", "SRC - Lecture Theatre 51"),
    @(47, "UTSRC-LT52", "SRC - Lecture Theatre 52"),
    @(48, "UTSRC-LT53", "SRC - Lecture Theatre 53"),
    @(49, "UTTP-AUD1", "TP- Auditorium 1"),
    @(50, "NA", "NA"),
    @(51, "ERC-ALR", "ERC - Active Learning Room"),
    @(52, "NAK-AUD", "ERC - Ngee Ann Kong Si Auditorium"),
    @(53, "Y-BlackBox", "ARTS CENTRE - BLACK BOX THEATRE"),
    @(54, "Y-ELMCL", "ELM COLLEGE - COMMON LOUNGE"),
    @(55, "Y-YESLab", "SCIENCE CENTRE - LAB 3 (Y.E.S. Lab)"),
    @(56, "TP-GLR", "TP - Global Learning Room")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2. Clear the leftover tall custom row height from the old layout (row
#    46 used to be a wrapped 24pt row; it is now a normal single-line
#    entry) by auto-fitting it back to the sheet's standard height.
# ---------------------------------------------------------------------
$ws.Rows.Item(46).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 3. Re-style + re-size the trailing legend block (rows 50-56): bump the
#    font to size 12 and set explicit row heights.
# ---------------------------------------------------------------------
$ws.Range("A50:B56").Font.Size = 12

$ws.Rows.Item(50).RowHeight = 14.25
$ws.Rows.Item(51).RowHeight = 16
$ws.Rows.Item(52).RowHeight = 16
$ws.Rows.Item(53).RowHeight = 16
$ws.Rows.Item(54).RowHeight = 16
$ws.Rows.Item(55).RowHeight = 16
$ws.Rows.Item(56).RowHeight = 16

# ---------------------------------------------------------------------
# 4. Fix up the hidden AutoFilter defined name so it matches the new,
#    shorter data extent (the filtered block used to run to row 53, now
#    it ends at row 49).
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Rooms_SD!_FilterDatabase") {
        $n.RefersTo = "=Rooms_SD!`$A`$31:`$B`$49"
    }
}

# ---------------------------------------------------------------------
# 5. Restore the selection the author left the sheet on.
# ---------------------------------------------------------------------
$ws.Range("B37").Select()
